$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: new work-log entry (Creation des Maquettes follow-up)
$ws.Cells.Item(17, 1).Value = 44686
$ws.Cells.Item(17, 2).Value = 0.5625
$ws.Cells.Item(17, 3).Value = 0.60416666666666663
$ws.Range("E17").Value = "Création Maquettes "

# Row 16: add the "Titre"/Description follow-up for the mockups task
$ws.Range("F16").Value = "Création Main Page, SignUp Page et Create Program"

$ws.Range("F17").Value = "Création Personal Programs, Personal Program, Exercice details, Create exercise"

# Row 18: MCD correction entry, adding the "programs" table note
$ws.Cells.Item(18, 1).Value = 44686
$ws.Cells.Item(18, 2).Value = 0.60416666666666663
$ws.Cells.Item(18, 3).Value = 0.62847222222222221
$ws.Range("E18").Value = "Correction MCD"
$ws.Range("F18").Value = "J'ai du ajouter une table ""programs"" pour les programmes des utilisateurs"

# Move the visible selection to match the saved view state
[void]$ws.Range("F19").Select()
